$wb = $excel.ActiveWorkbook

# ---- 1. Update the "总计" (summary) sheet: insert the new 2022-Q4 row and
#         shift the existing Q3/Q2/Q1 rows down by one ----
$summary = $wb.Worksheets.Item(1)

# Clone the formatting of the existing "index" column (A2:A4, style index 2)
# down onto the new A5 cell before we touch any values (same-sheet Copy keeps
# the original cellXf / shared style index instead of inventing a new one).
$summary.Cells.Item(4,1).Copy($summary.Cells.Item(5,1))

# Target state (after the edit) for rows 2..5, columns A..D
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 7
$summary.Cells.Item(2,4).Value = 0.78

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q3"
$summary.Cells.Item(3,3).Value = 13
$summary.Cells.Item(3,4).Value = 1.03

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2022-Q2"
$summary.Cells.Item(4,3).Value = 1
$summary.Cells.Item(4,4).Value = 0.68

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2022-Q1"
$summary.Cells.Item(5,3).Value = 2
$summary.Cells.Item(5,4).Value = 0.74

# ---- 2. Insert a brand-new "2022-Q4" sheet right before "2022-Q3" ----
# Cross-sheet Range.Copy doesn't transfer content/format in this host, so
# duplicate the whole "2022-Q3" sheet (Worksheet.Copy keeps every cellXf
# untouched) and overwrite its data in place instead.
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# 2022-Q3 has 13 data rows (rows 2-14); 2022-Q4 only needs 7 (rows 2-8)
$q4.Rows("9:14").Delete()

# header row text (already styled s="2" from the duplicated sheet)
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q4.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# data rows — column A (0..6 index) is already correct from the duplicated
# sheet, so only overwrite B..H
$rows = @(
    @("000242", "景顺长城策略精选A", "10.70", "93.00", "6.00", "0.6420", 1),
    @("001703", "银华沪港深增长股票A", "1.71", "93.53", "3.76", "0.0643", 7),
    @("016307", "景顺长城北交所精选两年定开混合A", "1.83", "43.56", "2.22", "0.0406", 7),
    @("001744", "诺安进取回报灵活配置混合", "0.59", "69.55", "3.55", "0.0209", 1),
    @("014364", "银华沪港深增长股票C", "0.25", "93.53", "3.76", "0.0094", 7),
    @("016308", "景顺长城北交所精选两年定开混合C", "0.27", "43.56", "2.22", "0.0060", 7),
    @("017167", "景顺长城策略精选C", "0.00", "93.00", "6.00", "0.0000", 1)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    $q4.Cells.Item($excelRow, 2).Value = "'" + $row[0]
    $q4.Cells.Item($excelRow, 3).Value = $row[1]
    $q4.Cells.Item($excelRow, 4).Value = "'" + $row[2]
    $q4.Cells.Item($excelRow, 5).Value = "'" + $row[3]
    $q4.Cells.Item($excelRow, 6).Value = "'" + $row[4]
    $q4.Cells.Item($excelRow, 7).Value = "'" + $row[5]
    $q4.Cells.Item($excelRow, 8).Value = $row[6]
}

# last data row (row 8) stores a plain numeric 0 in column G (not text "0.0000")
$q4.Cells.Item(8, 7).Value = 0
